$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pre-Alert Template Import")

$ws.Range("A3").Value = "JSSO1000249"
$ws.Range("B3").Value = "JSSO1000249"
$ws.Range("C3").Value = "JSSO1000249"
$ws.Range("AJ3").Value = "JSCN1000249"
$ws.Range("AL3").Value = "SLJSSO1000249"
$ws.Range("AN3").Value = "MBLJSSO1000249"
$ws.Range("AO3").Value = "HBLJSSO1000249"
